$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "Answers"
$ws.Range("C1").Value = "Source"
$ws.Range("D1").Value = "Source Link"
$ws.Range("E1").Value = "Dates Used"

$ws.Range("A2").Value = "Are there children in your household who are in Kindergarten through 12th grade?"
$ws.Range("B2").Value = "• Yes`n• No"
$ws.Range("C2").Value = "Developed by RAPID Team"
$ws.Range("E2").Value = "Aug 24 - Aug 27, Sept 8 - Sept 10"

$ws.Range("A3").Value = "[If yes to children in K-12]`n`nWhat has the school/school district your child(ren) plans to attend/currently attending decided regarding the school-year?"
$ws.Range("B3").Value = "• Open completely`n• Move to online only`n• A combination of in-person or some days and online on other days`n• Ability for parent to choose between online or in-person instruction`n• Other`n• Not applicable (e.g., homeschooled)`n• Don't know yet"
$ws.Range("C3").Value = "Developed by RAPID Team"
$ws.Range("E3").Value = "Aug 24 - Aug 27, Sept 8 - Sept 10"

$ws.Range("A4").Value = "[If yes to children in K-12]`n`nFor the time your child will be learning remotely, who will be assisting with your child's online learning? Select all that apply."
$ws.Range("B4").Value = "• You/other parent/step-parent`n• Sibling 15 years or older`n• Grandparent`n• Other relative`n• Friend of parent `n• Other neighbor`n• Adult responsible for a group (e.g., pod, bubble)`n• Babysitter/nanny/au pair`n• Other`n• No one is able to do this`n• Not applicable"
$ws.Range("C4").Value = "Developed by RAPID Team"
$ws.Range("E4").Value = "Aug 24 - Aug 27, Sept 8 - Sept 10"

$ws.Range("A5").Value = "[If yes to children in K-12]`n`nIf  your child has already begun the 2020-20Aug 24 - Aug 27 school year, how are they reacting to online learning?"
$ws.Range("B5").Value = "•`tVery well`n•`tSomewhat well`n•`tWell`n•`tNot well `n•`tVery poorly"
$ws.Range("C5").Value = "Developed by RAPID Team"
$ws.Range("E5").Value = "Sept 8 - Sept 10"

$ws.Range("A6").Value = "[If yes to children in K-12]`n`nIncluding hours spent during weekdays and weekends, about how many hours did you spend on teaching activities with your school-aged child(ren) in this household during the last 7 days? Enter the total number of hours. If none, enter 0."
$ws.Range("B6").Value = "Open Response"
$ws.Range("C6").Value = "RAPID Team Modified `nU.S. Census Bureau Household Pulse Survey "
$ws.Range("D6").Value = "https://www.census.gov/data/experimental-data-products/household-pulse-survey.html"
$ws.Range("E6").Value = "Sept 8 - Sept 10"

$ws.Range("A7").Value = "[If yes to children in K-12]`n`nWill you use a child care provider to help you when your child(ren) is not in school?"
$ws.Range("B7").Value = "• Yes`n• No`n• Maybe"
$ws.Range("C7").Value = "Developed by RAPID Team"
$ws.Range("E7").Value = "Aug 24 - Aug 27, Sept 8 - Sept 10"

$ws.Range("A8").Value = "If you have a child that was due to be entering kindergarten this fall, have you decided to wait until next fall (20Aug 24 - Aug 27) instead because of the pandemic?"
$ws.Range("B8").Value = "•`tYes`n•`tNo `n•`tNot applicable"
$ws.Range("C8").Value = "Developed by RAPID Team"
$ws.Range("E8").Value = "Sept 8 - Sept 10"

$ws.Range("A9").Value = "If yes, why? Select all that apply. "
$ws.Range("B9").Value = "•`tSafety`n•`tUncertain about the plan for school (in person/online)`n•`tNot able to manage online instruction for my child along with my other responsibilities (work, etc.)"
$ws.Range("C9").Value = "Developed by RAPID Team"
$ws.Range("E9").Value = "Sept 8 - Sept 10"

$ws.Range("A10").Value = "[If `"Will you use a child care provider to help you when your child(ren) is not in school?`" = Yes or Maybe]`n`nDoes this affect the child care arrangements you have for your child(ren) age 5 and under?"
$ws.Range("B10").Value = "• Yes`n• No`n• Maybe`n• Not applicable"
$ws.Range("C10").Value = "Developed by RAPID Team"
$ws.Range("E10").Value = "Aug 24 - Aug 27, Sept 8 - Sept 10"

$ws.Range("A11").Value = "For your child(ren) between the age of 0-5, which of the following are you doing to support their learning? Select all that apply."
$ws.Range("B11").Value = "• Using in-person games and activities at your home`n• Attending in-person learning outside of your home (daycare, childcare, etc.)`n• Using online resources (e.g., educational apps, educational TV shows, etc.)`n• Attending online classes/activities (facilitated by someone outside of your home, e.g., remote preschool, only story time, etc.)`n• Other`n• None of the above"
$ws.Range("C11").Value = "Developed by RAPID Team"
$ws.Range("E11").Value = "Aug 24 - Aug 27, Sept 8 - Sept 10"

$ws.Range("A12").Value = "Who will be assisting with your child(s) 0-5 learning? Select all that apply."
$ws.Range("B12").Value = "• You/other parent/step-parent`n• Sibling 15 years or older`n• Grandparent`n• Other relative`n• Friend of parent `n• Other neighbor`n• Adult responsible for a group (e.g., pod, bubble)`n• Babysitter/nanny/au pair`n• Other`n• No one is able to do this`n• Not applicable"
$ws.Range("C12").Value = "Developed by RAPID Team"
$ws.Range("E12").Value = "Aug 24 - Aug 27, Sept 8 - Sept 10"

$ws.Range("A13").Value = "Including hours spent during weekdays and weekends, about how many hours did you spend on face-to-face activities with your children 0-5 in this household during the last 7 days? `n`nEnter the total number of hours. If none, enter 0."
$ws.Range("B13").Value = "Open Response"
$ws.Range("C13").Value = "RAPID Team Modified `nU.S. Census Bureau Household Pulse Survey "
$ws.Range("D13").Value = "https://www.census.gov/data/experimental-data-products/household-pulse-survey.html"
$ws.Range("E13").Value = "Sept 8 - Sept 10"

$ws.Range("A14").Value = "Including hours spent during weekdays and weekends, about how many hours did you spend working during the last 7 days? `n`nEnter the total number of hours. If none, enter 0 or NA if not working currently. "
$ws.Range("B14").Value = "Open Response"
$ws.Range("C14").Value = "RAPID Team Modified `nU.S. Census Bureau Household Pulse Survey "
$ws.Range("D14").Value = "https://www.census.gov/data/experimental-data-products/household-pulse-survey.html"
$ws.Range("E14").Value = "Sept 8 - Sept 10"

$ws.Range("E1").Select() | Out-Null

Write-Output "Edit complete"